# Commit: Added dataset "res1C_focus_hand" to data-raw/
#
# This workbook (res1A_focus_hand.xlsx) was renamed/re-derived as part of a
# family of per-participant "focus_hand" coding workbooks. The sheet tab
# names get their hyphens swapped for underscores, and the last-touched
# cell on the second ("slg") sheet moves from I6 down to I25.

$wb = $excel.ActiveWorkbook

# --- Rename the two sheet tabs: hyphen -> underscore -------------------
$mdgSheet = $wb.Worksheets.Item(1)
$mdgSheet.Name = "mdg_20171117"

$slgSheet = $wb.Worksheets.Item(2)
$slgSheet.Name = "slg_20171117"

# --- Update the active selection on the "slg" sheet --------------------
# It is the tab-selected sheet, and its last recorded selection moves
# from I6 to I25.
$slgSheet.Activate()
$slgSheet.Range("I25").Select()
